# Re-insert the "+" continuation markers that were stripped from several
# footnote variant-reading notes, and fix two footnotes whose placeholder
# underscore ("_") should really be a plain space between the two "shad"
# (Tibetan punctuation marks, U+0F0D/ "|").
#
# The footnote pane text can only be read/written reliably through the
# Footnote's own Range.Text property in this runtime (Range.Find does not
# reach footnote story ranges), so we match each footnote by its current
# (old) text and overwrite it with the corrected (new) text.

$d = $word.ActiveDocument

# old footnote text -> new footnote text
$map = @{
    "མཐུན་པའི། སྣར་ཐང་། པེ་ཅིན།"            = "+མཐུན་པའི། སྣར་ཐང་། པེ་ཅིན།";
    "རང་གི་བཞིན། སྣར་ཐང་། པེ་ཅིན།"          = "+རང་གི་བཞིན། སྣར་ཐང་། པེ་ཅིན།";
    "བྱང་ཆུབ་ཀྱི། སྣར་ཐང་། པེ་ཅིན།"         = "+བྱང་ཆུབ་ཀྱི། སྣར་ཐང་། པེ་ཅིན།";
    "རྣམ་པར། པེ་ཅིན།"                       = "+རྣམ་པར། པེ་ཅིན།";
    "ཤེས་བྱ་བ། སྣར་ཐང་། པེ་ཅིན།"            = "+ཤེས་བྱ་བ། སྣར་ཐང་། པེ་ཅིན།";
    "བཀྲ་ཤིས་པར་གྱུར་ཅིག། སྣར་ཐང་།"         = "+བཀྲ་ཤིས་པར་གྱུར་ཅིག། སྣར་ཐང་།";
    "།།_། ཞེས་པར་མ་གཞན་ནང་མེད།"             = "།། ། ཞེས་པར་མ་གཞན་ནང་མེད།";
    "།_། ཞེས་པར་མ་གཞན་ནང་མེད།"              = "། ། ཞེས་པར་མ་གཞན་ནང་མེད།";
}

for ($i = 1; $i -le $d.Footnotes.Count; $i++) {
    $fn = $d.Footnotes.Item($i)
    $old = $fn.Range.Text
    if ($map.ContainsKey($old)) {
        $fn.Range.Text = $map[$old]
    }
}
